# Edit script: adds a "ground_truth" classification column (Q) to the sheet
# and removes two duplicate rows (old rows 13 & 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two duplicate rows (old rows 13 and 14 are exact duplicates
#    of row 12). Everything below shifts up automatically.
$ws.Rows.Item(13).Resize(2).Delete()

# 2) Add the new "ground_truth" header in Q1, re-using the same bold font
#    as the other header cells (A1:P1) but with a left/right-only border.
$ws.Range("A1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("Q1").Borders.Item(8).LineStyle = -4142
$ws.Range("Q1").Borders.Item(9).LineStyle = -4142
$ws.Range("Q1").Value = "ground_truth"

# 3) Fill in the classification values for the data rows (Q2:Q22).
$ws.Cells.Item(2, 17).Value = "Hilfsreport"
$ws.Cells.Item(3, 17).Value = "Andere"
$ws.Cells.Item(4, 17).Value = "Andere"
$ws.Cells.Item(5, 17).Value = "Andere"
$ws.Cells.Item(6, 17).Value = "Hilfsreport"
$ws.Cells.Item(7, 17).Value = "Hilfsreport"
$ws.Cells.Item(8, 17).Value = "Hilfsreport"
$ws.Cells.Item(9, 17).Value = "Hilfsreport"
$ws.Cells.Item(10, 17).Value = "Andere"
$ws.Cells.Item(11, 17).Value = "Andere"
$ws.Cells.Item(12, 17).Value = "Hilfsreport"
$ws.Cells.Item(13, 17).Value = "Hilfsreport"
$ws.Cells.Item(14, 17).Value = "Hilfsreport"
$ws.Cells.Item(15, 17).Value = "Hilfsreport"
$ws.Cells.Item(16, 17).Value = "Andere"
$ws.Cells.Item(17, 17).Value = "Andere"
$ws.Cells.Item(18, 17).Value = "Andere"
$ws.Cells.Item(19, 17).Value = "Andere"
$ws.Cells.Item(20, 17).Value = "Andere"
$ws.Cells.Item(21, 17).Value = "Andere"
$ws.Cells.Item(22, 17).Value = "Andere"

# 4) Restore the view/selection state.
$ws.Range("P12").Select()
